# daily auto push: 2026-02-16 10:06 UTC
#
# Insert a new daily-log row before row 831 (this shifts the existing
# 2026/12/29 .. 2027/01/05 rows down by one, growing the used range from
# A1:D872 to A1:D873), then populate the newly inserted row with the
# day's data: 2026/02/16 (Mon), hour 17, value 201.
#
# Column A holds the date as literal text (matching every other row in
# the sheet, which stores dates as plain strings, not Excel date
# serials). Assigning a slash-delimited string straight to .Value would
# get auto-recognized as a date by Excel, so the cell is pre-formatted
# as Text, the literal is written, and the formatting is cleared back
# to the sheet's default afterwards (leaving the cell's style exactly
# like its neighbors).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(831).Insert()

$ws.Cells.Item(831, 1).NumberFormat = "@"
$ws.Cells.Item(831, 1).Value = "2026/02/16"
$ws.Cells.Item(831, 1).ClearFormats()

$ws.Cells.Item(831, 2).Value = "月"
$ws.Cells.Item(831, 3).Value = 17
$ws.Cells.Item(831, 4).Value = 201
